$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.797.82'
$ws.Range('E2').Value = '  +4.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.760.02'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.52'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.55'
$ws.Range('E6').Value = '  +6.47%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +2.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.757.58'
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('E11').Value = '  +5.30%  '
$ws.Range('E12').Value = '  +4.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.389'
$ws.Range('E13').Value = '  +3.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.240.53'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.44'
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.744.32'
$ws.Range('E16').Value = '  +4.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000153'
$ws.Range('E17').Value = '  +6.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.754.63'
$ws.Range('E18').Value = '  +3.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.98'
$ws.Range('E19').Value = '  +3.65%  '
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '361.35'
$ws.Range('E21').Value = '  +3.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.99'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.14'
$ws.Range('E25').Value = '  +3.26%  '
$ws.Range('E26').Value = '  +4.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.56'
$ws.Range('E27').Value = '  +4.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0915'
$ws.Range('E29').Value = '  +12.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.00'
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.15'
$ws.Range('E31').Value = '  +6.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '171.57'
$ws.Range('E32').Value = '  +2.73%  '
$ws.Range('E33').Value = '  +15.07%  '
$ws.Range('E35').Value = '  +3.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.79'
$ws.Range('E36').Value = '  +7.77%  '
$ws.Range('E37').Value = '  +9.24%  '
$ws.Range('E38').Value = '  +10.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.01'
$ws.Range('E39').Value = '  +14.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '345.15'
$ws.Range('E40').Value = '  +5.35%  '
$ws.Range('E41').Value = '  +5.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.41'
$ws.Range('E42').Value = '  +2.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.62'
$ws.Range('E43').Value = '  +7.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.74'
$ws.Range('E44').Value = '  +5.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.83'
$ws.Range('E45').Value = '  +6.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '140.09'
$ws.Range('E46').Value = '  +3.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0591'
$ws.Range('E47').Value = '  +5.71%  '
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('E49').Value = '  +3.64%  '
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('E51').Value = '  -0.24%  '
